# Applies the commit "Elimna EC anteriores y se agregan nuevos, se modifica
# base de datos" to the account-statement worksheet:
#   - refresh the detail table (rows 36-53) with the new arrears rows
#   - move the signature/footer block down to stay right after the table
#   - refresh the summary totals (Valor Mora / Cant. Trabajadores / Cant. Periodos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 9 blank rows before row 45 to make room for the extra
#     detail rows. This shifts the footer block (old rows 49-50) down to
#     rows 58-59 and keeps merged cells / the sheet dimension in sync. ---
$ws.Rows("45:53").Insert(-4121)

# --- Step 2: row 44 used to be the last row of the table, so it carries the
#     special "closing" bottom-border style. Copy that style onto the new
#     last row (53) first, then restyle row 44 itself like a normal interior
#     row (matching row 43) because more rows now follow it. ---
$ws.Range("B44:J44").Copy()
$ws.Range("B53:J53").PasteSpecial(-4122)
$ws.Range("B43:J43").Copy()
$ws.Range("B44:J44").PasteSpecial(-4122)

# --- Step 3: stamp normal interior-row formatting onto the rest of the
#     freshly inserted rows (45-52). ---
$ws.Range("B16:J16").Copy()
$ws.Range("B45:J52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: write the refreshed detail data (rows 36-53). ---
$data = @(
  @(36, "PE", "834804928091987", "RAFAEL ANTONIO MARIN GUERRERO", "2311", 46400, 1160000),
  @(37, "PE", "834804928091987", "RAFAEL ANTONIO MARIN GUERRERO", "2310", 46400, 1160000),
  @(38, "PE", "834804928091987", "RAFAEL ANTONIO MARIN GUERRERO", "2309", 46400, 1160000),
  @(39, "PE", "834804928091987", "RAFAEL ANTONIO MARIN GUERRERO", "2308", 46400, 1160000),
  @(40, "CC", "1001832933", "PIERINA SERGE PEREZ", "2507", 40000, 1000000),
  @(41, "CC", "1001832933", "PIERINA SERGE PEREZ", "2506", 40000, 1000000),
  @(42, "CC", "1001832933", "PIERINA SERGE PEREZ", "2505", 40000, 1000000),
  @(43, "CC", "1001832933", "PIERINA SERGE PEREZ", "2504", 40000, 1000000),
  @(44, "CC", "1001832933", "PIERINA SERGE PEREZ", "2503", 40000, 1000000),
  @(45, "CC", "1001832933", "PIERINA SERGE PEREZ", "2502", 40000, 1000000),
  @(46, "CC", "1002189028", "ANGELLYS SARAIS PIANETA JULIO", "2507", 56940, 1423500),
  @(47, "CC", "1002189028", "ANGELLYS SARAIS PIANETA JULIO", "2506", 56940, 1423500),
  @(48, "CC", "1002189028", "ANGELLYS SARAIS PIANETA JULIO", "2505", 56940, 1423500),
  @(49, "CC", "1049827991", "DIOMARA DIAZ HERRERA", "2507", 56940, 1423500),
  @(50, "CC", "1049827991", "DIOMARA DIAZ HERRERA", "2506", 56940, 1423500),
  @(51, "CC", "1049827991", "DIOMARA DIAZ HERRERA", "2505", 56940, 1423500),
  @(52, "CC", "1049827991", "DIOMARA DIAZ HERRERA", "2504", 56940, 1423500),
  @(53, "CC", "1049827991", "DIOMARA DIAZ HERRERA", "2503", 56940, 1423500)
)
foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}

# --- Step 5: refresh the summary box with the new totals. ---
$ws.Range("E11").Value = 1809120   # Valor Mora total
$ws.Range("C13").Value = 4         # Cant. Trabajadores
$ws.Range("F13").Value = 24        # Cant. Periodos
